# Streamlit_BC.pptx update:
#  1. Bump the cached "datetimeFigureOut" field text (11/13/22 -> 11/14/22)
#     on the slide master and every slide layout - this mirrors the
#     automatic re-cache PowerPoint performs on the "today" date field.
#  2. Add a new textbox with the YouTube link to the last slide (33),
#     placed after the existing group shape.

$p = $ppt.ActivePresentation

$oldDate = "11/13/22"
$newDate = "11/14/22"

# --- 1. Update the date placeholder on the slide master ---
$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Update the date placeholder on every slide layout ---
for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $lay = $m.CustomLayouts.Item($j)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2. Add the YouTube-link textbox to the last slide ---
$s = $p.Slides.Item($p.Slides.Count)

# Burn through two auto-assigned shape ids so the real textbox below
# lands on id=5 (matching the authored file), mirroring how ids that
# were previously used/recycled in the deck get reassigned.
$burn1 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$burn1.Delete()
$burn2 = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$burn2.Delete()

# PowerPoint's Shapes.AddTextbox positions/sizes are expressed in points;
# the target geometry below is given in EMU (1 pt = 12700 EMU).
$emuPerPt = 12700
$left   = 6111240 / $emuPerPt
$top    = 3352800 / $emuPerPt
$width  = 4872616 / $emuPerPt
$height = 369332 / $emuPerPt

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)

$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = 0

$tr = $tb.TextFrame.TextRange
$tr.Text = "https://"
$tr = $tr.InsertAfter("www.youtube.com")
$tr = $tr.InsertAfter("/")
$tr = $tr.InsertAfter("watch?v")
$tr = $tr.InsertAfter("=QCucIBfCy84")
